$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.558.99"
$ws.Range("E2").Value = "  +0.46%  "

# Row 3
$ws.Range("D3").Value = "1.818.53"
$ws.Range("E3").Value = "  +0.98%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "226.00"
$ws.Range("E5").Value = "  +0.37%  "

# Row 6
$ws.Range("D6").Value = "0.606"
$ws.Range("E6").Value = "  +1.16%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "43.83"
$ws.Range("E8").Value = "  +20.97%  "

# Row 9
$ws.Range("D9").Value = "0.295"
$ws.Range("E9").Value = "  +1.15%  "

# Row 10
$ws.Range("D10").Value = "0.0676"
$ws.Range("E10").Value = "  -0.15%  "

# Row 11
$ws.Range("D11").Value = "0.100"
$ws.Range("E11").Value = "  +3.72%  "

# Row 12
$ws.Range("D12").Value = "2.079.90"
$ws.Range("E12").Value = "  +0.91%  "

# Row 13
$ws.Range("D13").Value = "1.817.78"
$ws.Range("E13").Value = "  +0.72%  "

# Row 14
$ws.Range("D14").Value = "11.14"
$ws.Range("E14").Value = "  -1.19%  "

# Row 15
$ws.Range("D15").Value = "0.639"
$ws.Range("E15").Value = "  +1.82%  "

# Row 16
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "4.48"
$ws.Range("E16").Value = "  +1.02%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "34.532.57"
$ws.Range("E17").Value = "  +0.42%  "

# Row 18
$ws.Range("D18").Value = "67.86"
$ws.Range("E18").Value = "  -0.93%  "

# Row 19
$ws.Range("D19").Value = "242.23"
$ws.Range("E19").Value = "  -0.02%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0780"
$ws.Range("E20").Value = "  +1.00%  "

# Row 21
$ws.Range("D21").Value = "11.57"
$ws.Range("E21").Value = "  +2.85%  "

# Row 22
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("D23").Value = "4.52"
$ws.Range("E23").Value = "  +10.52%  "

# Row 24
$ws.Range("E24").Value = "  -2.19%  "

# Row 25
$ws.Range("D25").Value = "171.42"
$ws.Range("E25").Value = "  +0.42%  "

# Row 26
$ws.Range("D26").Value = "7.79"
$ws.Range("E26").Value = "  -1.13%  "

# Row 27
$ws.Range("D27").Value = "17.56"
$ws.Range("E27").Value = "  +1.21%  "

# Row 28
$ws.Range("E28").Value = "  +0.59%  "

# Row 29
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("D30").Value = "3.85"
$ws.Range("E30").Value = "  +1.62%  "

# Row 31
$ws.Range("E31").Value = "  +0.64%  "

# Row 32
$ws.Range("D32").Value = "3.92"
$ws.Range("E32").Value = "  +0.27%  "

# Row 33
$ws.Range("D33").Value = "0.0519"
$ws.Range("E33").Value = "  +0.91%  "

# Row 34
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +2.99%  "

# Row 35
$ws.Range("D35").Value = "89.64"
$ws.Range("E35").Value = "  +11.11%  "

# Row 36
$ws.Range("D36").Value = "0.659"
$ws.Range("E36").Value = "  +1.45%  "

# Row 37
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.323.43"
$ws.Range("E37").Value = "  -2.83%  "

# Row 38
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "15.21"
$ws.Range("E38").Value = "  +14.61%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "2.41"
$ws.Range("E39").Value = "  +2.09%  "

# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.06"
$ws.Range("E40").Value = "  +0.15%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0190"
$ws.Range("E41").Value = "  +2.61%  "

# Row 42
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.958"
$ws.Range("E42").Value = "  +2.25%  "

# Row 43
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "1.22"
$ws.Range("E43").Value = "  +4.77%  "

# Row 44
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  +1.41%  "

# Row 45
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "2.42"
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").Value = "0.0520"
$ws.Range("E46").Value = "  +4.10%  "

# Row 47
$ws.Range("D47").Value = "1.976.95"
$ws.Range("E47").Value = "  +0.71%  "

# Row 48
$ws.Range("D48").Value = "5.88"
$ws.Range("E48").Value = "  +1.70%  "

# Row 49
$ws.Range("E49").Value = "  +0.03%  "

# Row 50
$ws.Range("D50").Value = "101.45"
$ws.Range("E50").Value = "  -0.80%  "

# Row 51
$ws.Range("D51").Value = "0.0612"
$ws.Range("E51").Value = "  +1.41%  "
